$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PAS-730: the four "control table" symbol columns (BI_SYMBOL, PD_SYMBOL,
# UM_SYMBOL, MP_SYMBOL -> columns AE:AH) on row 6 all get collapsed onto the
# same code "X" to remove the date overlap with the other VOLKSWAGEN GOLF
# control rows. Previously these held different individual codes (R/E/S/A).
$ws.Range("AE6:AH6").Value = "X"

# Leave the sheet scrolled/selected on the cell that was last edited (AH6),
# matching the workbook's saved view state after the edit.
$ws.Range("AH6").Select()
